{"js": "// Replace the ten-problem \"two-digit x two-digit\" answer grid with a newly\n// generated set of problems/answers. Each old equation text is unique in\n// the document, so a straightforward search-and-replace per pair is safe.\nconst replacements = [\n  { old: \"40\u00d791=3640\", new: \"13\u00d718=234\" },\n  { old: \"81\u00d712=972\", new: \"30\u00d794=2820\" },\n  { old: \"69\u00d742=2898\", new: \"52\u00d734=1768\" },\n  { old: \"85\u00d730=2550\", new: \"21\u00d780=1680\" },\n  { old: \"42\u00d765=2730\", new: \"79\u00d753=4187\" },\n  { old: \"20\u00d742=840\", new: \"32\u00d742=1344\" },\n  { old: \"65\u00d759=3835\", new: \"58\u00d750=2900\" },\n  { old: \"93\u00d785=7905\", new: \"35\u00d787=3045\" },\n  { old: \"86\u00d747=4042\", new: \"86\u00d741=3526\" },\n  { old: \"18\u00d764=1152\", new: \"35\u00d765=2275\" },\n  { old: \"63\u00d795=5985\", new: \"73\u00d729=2117\" },\n  { old: \"70\u00d787=6090\", new: \"12\u00d798=1176\" },\n  { old: \"14\u00d751=714\", new: \"53\u00d727=1431\" },\n  { old: \"85\u00d763=5355\", new: \"40\u00d717=680\" },\n  { old: \"94\u00d784=7896\", new: \"50\u00d796=4800\" },\n  { old: \"53\u00d765=3445\", new: \"69\u00d740=2760\" },\n  { old: \"88\u00d780=7040\", new: \"44\u00d785=3740\" },\n  { old: \"50\u00d777=3850\", new: \"88\u00d740=3520\" },\n  { old: \"14\u00d784=1176\", new: \"47\u00d771=3337\" },\n  { old: \"83\u00d744=3652\", new: \"23\u00d790=2070\" },\n  { old: \"36\u00d774=2664\", new: \"97\u00d785=8245\" },\n  { old: \"23\u00d773=1679\", new: \"36\u00d773=2628\" },\n  { old: \"42\u00d753=2226\", new: \"93\u00d753=4929\" },\n  { old: \"48\u00d749=2352\", new: \"43\u00d721=903\" },\n  { old: \"69\u00d720=1380\", new: \"64\u00d791=5824\" },\n];\n\nfor (const { old, new: replacement } of replacements) {\n  const results = context.document.body.search(old, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the ten-problem \"two-digit x two-digit\" answer grid with a newly\n# generated set of problems/answers. Each old equation text is unique in\n# the document, so a straightforward Find/Replace per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"40\u00d791=3640\"; New = \"13\u00d718=234\" }\n    @{ Old = \"81\u00d712=972\"; New = \"30\u00d794=2820\" }\n    @{ Old = \"69\u00d742=2898\"; New = \"52\u00d734=1768\" }\n    @{ Old = \"85\u00d730=2550\"; New = \"21\u00d780=1680\" }\n    @{ Old = \"42\u00d765=2730\"; New = \"79\u00d753=4187\" }\n    @{ Old = \"20\u00d742=840\"; New = \"32\u00d742=1344\" }\n    @{ Old = \"65\u00d759=3835\"; New = \"58\u00d750=2900\" }\n    @{ Old = \"93\u00d785=7905\"; New = \"35\u00d787=3045\" }\n    @{ Old = \"86\u00d747=4042\"; New = \"86\u00d741=3526\" }\n    @{ Old = \"18\u00d764=1152\"; New = \"35\u00d765=2275\" }\n    @{ Old = \"63\u00d795=5985\"; New = \"73\u00d729=2117\" }\n    @{ Old = \"70\u00d787=6090\"; New = \"12\u00d798=1176\" }\n    @{ Old = \"14\u00d751=714\"; New = \"53\u00d727=1431\" }\n    @{ Old = \"85\u00d763=5355\"; New = \"40\u00d717=680\" }\n    @{ Old = \"94\u00d784=7896\"; New = \"50\u00d796=4800\" }\n    @{ Old = \"53\u00d765=3445\"; New = \"69\u00d740=2760\" }\n    @{ Old = \"88\u00d780=7040\"; New = \"44\u00d785=3740\" }\n    @{ Old = \"50\u00d777=3850\"; New = \"88\u00d740=3520\" }\n    @{ Old = \"14\u00d784=1176\"; New = \"47\u00d771=3337\" }\n    @{ Old = \"83\u00d744=3652\"; New = \"23\u00d790=2070\" }\n    @{ Old = \"36\u00d774=2664\"; New = \"97\u00d785=8245\" }\n    @{ Old = \"23\u00d773=1679\"; New = \"36\u00d773=2628\" }\n    @{ Old = \"42\u00d753=2226\"; New = \"93\u00d753=4929\" }\n    @{ Old = \"48\u00d749=2352\"; New = \"43\u00d721=903\" }\n    @{ Old = \"69\u00d720=1380\"; New = \"64\u00d791=5824\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
